# Revised errors in data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the Door Count values that were mistakenly entered as 4 instead of 2
$ws.Range("F11").Value = 2
$ws.Range("F18").Value = 2
$ws.Range("F25").Value = 2

# Update the saved view/selection state of the sheet: scroll back to top
# and move the active selection to F6
$ws.Range("F6").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
